# Add files via upload
# Populates "Week 6" timesheet rows 2-12 with new entries, which also
# ripples the Weekly/Project Total formulas forward through the rest of
# the workbook, and updates the active-sheet/selection bookkeeping.

$wb = $excel.ActiveWorkbook

$wsWeek5 = $wb.Worksheets.Item("Week 5")
$wsWeek6 = $wb.Worksheets.Item("Week 6")

# ---------------------------------------------------------------------
# New timesheet entries for Week 6 (rows 2-12)
# Columns: A=Date, B=Start Time, C=End Time, D=Description, E=Hours
# ---------------------------------------------------------------------

$rows = @(
    @{ Row=2;  Date=42050; Start=0.041666666666666664; End=0.08333333333333333;  Desc="BeerCityMaps Database Design";         Hours=1 }
    @{ Row=3;  Date=42051; Start=$null;                 End=$null;               Desc=$null;                                  Hours=$null }
    @{ Row=4;  Date=42052; Start=0.33333333333333331;   End=0.41666666666666669;  Desc="BeerCityMaps Database Design";         Hours=2 }
    @{ Row=5;  Date=42053; Start=0.29166666666666669;   End=0.375;                Desc="BeerCityMaps Database Design";         Hours=2 }
    @{ Row=6;  Date=42053; Start=0.41666666666666669;   End=0.5;                  Desc="Class";                                Hours=2 }
    @{ Row=7;  Date=42054; Start=0.33333333333333331;   End=0.41666666666666669;  Desc="BeerCityMaps database-connection.php"; Hours=2 }
    @{ Row=8;  Date=42055; Start=0.33333333333333331;   End=0.375;                Desc="BeerCityMaps database-connection.php"; Hours=1 }
    @{ Row=9;  Date=42055; Start=0.41666666666666669;   End=0.5;                  Desc="Class";                                Hours=2 }
    @{ Row=10; Date=42056; Start=0.35416666666666669;   End=0.4375;               Desc="BeerCityMaps Database Design";         Hours=2 }
    @{ Row=11; Date=42056; Start=0.47916666666666669;   End=0.14583333333333334;  Desc="BeerCityMaps Registration Page";       Hours=4 }
    @{ Row=12; Date=42056; Start=0.1875;                End=0.3125;               Desc="BeerCityMaps Registration Page";       Hours=3 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $wsWeek6.Range("A$rowNum").Value = $r.Date
    $wsWeek6.Range("A$rowNum").NumberFormat = "d-mmm"

    if ($null -ne $r.Start) {
        $wsWeek6.Range("B$rowNum").Value = $r.Start
        $wsWeek6.Range("B$rowNum").NumberFormat = "h:mm"
    }

    if ($null -ne $r.End) {
        $wsWeek6.Range("C$rowNum").Value = $r.End
        $wsWeek6.Range("C$rowNum").NumberFormat = "h:mm"
    }

    if ($null -ne $r.Desc) {
        $wsWeek6.Range("D$rowNum").Value = $r.Desc
    }

    if ($null -ne $r.Hours) {
        $wsWeek6.Range("E$rowNum").Value = $r.Hours
    }

    $wsWeek6.Rows.Item($rowNum).RowHeight = 18
}

# ---------------------------------------------------------------------
# Selection / active-tab bookkeeping: the edit session left off on
# Week 5 (cell D3) and most-recently touched Week 6 (cell E13), which
# becomes the active tab.
# ---------------------------------------------------------------------

$wsWeek5.Activate()
$wsWeek5.Range("D3").Select()

$wsWeek6.Activate()
$wsWeek6.Range("E13").Select()
